$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "abc"
$ws.Range("B2").Value = "def"

$ws.Range("B3").Select()
